$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status changes: "In Translation" -> "Ready for handoff"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date/Datetime updates
$wsOverview.Range("D2").Value = "2016-35-13 02:35:12"
$wsZhCn.Range("E2").Value = "2016-03-13 02:35:09"
$wsDeDe.Range("E2").Value = "2016-03-13 02:35:12"
